$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '95.286.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.605.14'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.52%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.13'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '653.31'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.50%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.37%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.02%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.00'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.605.19'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.54%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.10'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.42'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.275.97'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '95.221.17'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.602.70'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.83'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.00%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.90'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '507.88'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.477'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.00%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.11%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.42'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.799.65'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.47'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.07%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.22'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.139'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.47'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +9.29%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.560'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.07'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '558.27'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.93%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.912'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '36.35'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +43.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.75'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.67'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.56'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.27'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0412'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.59'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.31'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.20%  '
